$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.235.87"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "2.425.36"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0801"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").Value = "2.804.45"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "2.429.21"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.834"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").Value = "44.140.73"
$ws.Range("E18").Value = "  +2.43%  "
$ws.Range("E19").Value = "  +1.28%  "
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  +2.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  +1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("E29").Value = "  +4.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.00%  "
$ws.Range("E31").Value = "  +12.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.56%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0760"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.40%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "131.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +24.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.30"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("E41").Value = "  +0.94%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.64%  "
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("D44").Value = "1.949.83"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  +9.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("E51").Value = "  +0.99%  "
